$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row just above the current row 131 ("Primera" quality,
# Espinaca / Vega Modelo de Temuco row) for the latest weekly observation.
# This shifts the existing rows 131..264 down to 132..265.
$ws.Rows.Item(131).Insert()

$ws.Cells.Item(131, 1).Value = 10
$ws.Cells.Item(131, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(131, 3).Value = "La Araucanía"
$ws.Cells.Item(131, 4).Value = 45033
$ws.Cells.Item(131, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(131, 5).Value = 9
$ws.Cells.Item(131, 6).Value = 100112012
$ws.Cells.Item(131, 7).Value = "Espinaca"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 125
$ws.Cells.Item(131, 11).Value = 10000
$ws.Cells.Item(131, 12).Value = 10000
$ws.Cells.Item(131, 13).Value = 10000
$ws.Cells.Item(131, 14).Value = "`$/docena de atados"
$ws.Cells.Item(131, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(131, 16).Value = 3333
$ws.Cells.Item(131, 17).Value = 3
$ws.Cells.Item(131, 18).Value = "Hortaliza"
